# Update cryptocurrency price/volume data per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.863.58"
$ws.Range("E2").Value = "  +1.17%  "
$ws.Range("D3").Value = "2.087.89"
$ws.Range("E3").Value = "  +0.92%  "
$ws.Range("E4").Value = "  +0.00%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "235.16"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -0.08%  "
$ws.Range("E6").Value = "  -0.39%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "59.54"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +3.56%  "
$ws.Range("E8").Value = "  -0.01%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.393"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -0.80%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.0793"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +2.60%  "
$ws.Range("E11").Value = "  +3.12%  "
$ws.Range("D12").Value = "2.394.59"
$ws.Range("E12").Value = "  +0.91%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "14.75"
$c.Style = "Normal"
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "21.47"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +3.44%  "
$ws.Range("E15").Value = "  -1.32%  "
$ws.Range("E16").Value = "  +2.28%  "
$ws.Range("D17").Value = "2.084.82"
$ws.Range("E17").Value = "  +0.78%  "
$ws.Range("D18").Value = "37.777.30"
$ws.Range("E18").Value = "  +1.17%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "6.24"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -1.80%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "71.86"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +3.06%  "
$ws.Range("D21").Value = "0.0₃0831"
$ws.Range("E21").Value = "  +1.62%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "229.14"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +1.10%  "
$ws.Range("E23").Value = "  +0.03%  "
$ws.Range("E24").Value = "  -0.60%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "2.41"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +0.37%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "170.94"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +2.20%  "
$ws.Range("E27").Value = "  +9.31%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "9.10"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +2.51%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "1.43"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +0.24%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "19.61"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +2.56%  "
$ws.Range("E31").Value = "  +1.65%  "
$ws.Range("E32").Value = "  +4.32%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "0.0634"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +2.37%  "
$ws.Range("E34").Value = "  +3.13%  "
$ws.Range("E35").Value = "  +0.70%  "
$ws.Range("E36").Value = "  +6.10%  "
$ws.Range("E37").Value = "  +2.74%  "
$ws.Range("E38").Value = "  -0.16%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "5.48"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -3.39%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.0987"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +2.24%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "99.73"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +1.40%  "
$ws.Range("E42").Value = "  -0.01%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.0216"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +1.35%  "
$ws.Range("D44").Value = "1.465.37"
$ws.Range("E44").Value = "  -1.12%  "
$ws.Range("E45").Value = "  +0.74%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "4.33"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +6.88%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "16.22"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +5.53%  "
$ws.Range("E48").Value = "  +4.32%  "
$ws.Range("E49").Value = "  +3.10%  "
$ws.Range("E50").Value = "  +2.57%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "47.57"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +5.86%  "
